$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.643.41"
$ws.Range("E2").Value = "  -1.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.886.19"
$ws.Range("E3").Value = "  -1.77%  "

$ws.Range("E4").Value = "  +0.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.02"
$ws.Range("E5").Value = "  -4.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4874"
$ws.Range("E7").Value = "  -2.69%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2885"
$ws.Range("E8").Value = "  -4.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06653"
$ws.Range("E9").Value = "  -3.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.881.02"
$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.69"
$ws.Range("E11").Value = "  -2.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07231"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "88.78"
$ws.Range("E13").Value = "  -2.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.001"
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6633"
$ws.Range("E15").Value = "  -3.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.595.66"
$ws.Range("E16").Value = "  -1.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007845"
$ws.Range("E17").Value = "  -2.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  +0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  -3.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.124.50"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.727"
$ws.Range("E22").Value = "  -3.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "187.65"
$ws.Range("E23").Value = "  +0.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.035"
$ws.Range("E24").Value = "  -1.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.259"
$ws.Range("E25").Value = "  -1.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.21"
$ws.Range("E26").Value = "  +3.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.24"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.828"
$ws.Range("E28").Value = "  -6.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.403"
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.245"
$ws.Range("E30").Value = "  -2.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09024"
$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.926"
$ws.Range("E32").Value = "  -3.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05176"
$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7303"
$ws.Range("E34").Value = "  -3.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.080"
$ws.Range("E35").Value = "  -6.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.691"
$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01821"
$ws.Range("E37").Value = "  -5.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.660"
$ws.Range("E38").Value = "  -3.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9182"
$ws.Range("E39").Value = "  -2.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.044"
$ws.Range("E40").Value = "  -7.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4318"
$ws.Range("E41").Value = "  -1.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.80"
$ws.Range("E42").Value = "  -1.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9988"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.723"
$ws.Range("E44").Value = "  -3.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1340"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.259"
$ws.Range("E46").Value = "  -8.43%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4013"
$ws.Range("E47").Value = "  +2.45%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05825"
$ws.Range("E48").Value = "  -0.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.625"
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.403"
$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.14"
$ws.Range("E51").Value = "  -0.86%  "
